# Apply the PNAD 2009 "agressao" correction:
#  - B2 header changes from "unnamed: 1_level_1" to "total"
#  - The two label-only rows ("situação do domicílio" / row 5 and
#    "grandes regiões" / row 8) are removed, so every data row below
#    shifts up to close the gaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled header in row 2.
$ws.Range("B2").Value = "total"

# Remove the two empty "section label" rows. Deleting the higher-numbered
# row first keeps the lower row's index ("row 5") valid for the second
# delete.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
